$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column header "M" in AJ1 (new shared string, index 35)
$ws.Range("AJ1").Value = "M"

# Row 2: first formula entered standalone (becomes the plain, non-shared formula)
$ws.Range("AJ2").Formula = "=SUM(F2:AI2)"

# Rows 3-61: formula filled down as one block (becomes a shared formula group
# anchored at AJ3, matching how the author extended the formula after typing
# it once in AJ2 and then filling/copying it down the rest of the column)
$ws.Range("AJ3:AJ61").Formula = "=SUM(F3:AI3)"

# Leave the selection on AJ6, matching the cursor position left by the author
$ws.Range("AJ6").Select()
